$wb = $excel.ActiveWorkbook

# --- Sheet '展览' (Exhibitions) ---
$ws1 = $wb.Worksheets.Item('展览')
# row 2
$ws1.Range('B2').NumberFormat = '@'
$ws1.Range('B2').Value = '2024-10-06'
$ws1.Range('B2').Style = 'Normal'
$ws1.Range('C2').Value = '南昌·星辰动漫游戏展嘉年华'
$ws1.Range('D2').Value = '长寿路 锐成体育2＋1篮球公园'
$ws1.Range('E2').Value = '2024.10.06 10:00-10.06 17:00'
$ws1.Range('F2').Value = 165
$ws1.Range('G2').Value = 45
$ws1.Range('H2').Value = 'https://show.bilibili.com/platform/detail.html?id=91961'
$ws1.Range('I2').Value = '//i2.hdslb.com/bfs/openplatform/202409/A6PiZvto1725362505262.jpeg'
# row 3
$ws1.Range('B3').NumberFormat = '@'
$ws1.Range('B3').Value = '2024-10-19'
$ws1.Range('B3').Style = 'Normal'
$ws1.Range('C3').Value = '南昌·第一届哥布林动漫游戏展——开学季&贺中秋'
$ws1.Range('D3').Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws1.Range('E3').Value = '2024.10.19 10:00-10.20 18:00'
$ws1.Range('F3').Value = 648
$ws1.Range('G3').Value = 55
$ws1.Range('H3').Value = 'https://show.bilibili.com/platform/detail.html?id=89240'
$ws1.Range('I3').Value = '//i2.hdslb.com/bfs/openplatform/202409/3N3iIqVu1725270119618.jpeg'
# row 4
$ws1.Range('B4').NumberFormat = '@'
$ws1.Range('B4').Value = '2024-11-02'
$ws1.Range('B4').Style = 'Normal'
$ws1.Range('C4').Value = '南昌·花绒万兽秋镜派对'
$ws1.Range('D4').Value = '双港西大街899号 旭辉Cmall(南昌店)'
$ws1.Range('E4').Value = '2024.11.02 10:00-11.03 21:30'
$ws1.Range('F4').Value = 24
$ws1.Range('G4').Value = 168
$ws1.Range('H4').Value = 'https://show.bilibili.com/platform/detail.html?id=92859'
$ws1.Range('I4').Value = '//i2.hdslb.com/bfs/openplatform/202409/7hJL2m3F1727175584690.jpeg'
# row 5
$ws1.Range('B5').NumberFormat = '@'
$ws1.Range('B5').Value = '2024-11-16'
$ws1.Range('B5').Style = 'Normal'
$ws1.Range('C5').Value = '上饶·星河城市动漫文化节'
$ws1.Range('D5').Value = '春江北大道时光PARK内 博悦宴会艺术中心'
$ws1.Range('E5').Value = '2024.11.16 10:00-11.16 17:00'
$ws1.Range('F5').Value = 214
$ws1.Range('G5').Value = 55
$ws1.Range('H5').Value = 'https://show.bilibili.com/platform/detail.html?id=92572'
$ws1.Range('I5').Value = '//i2.hdslb.com/bfs/openplatform/202409/xp4jNVRG1727165677359.jpeg'
# row 6
$ws1.Range('B6').NumberFormat = '@'
$ws1.Range('B6').Value = '2024-11-16'
$ws1.Range('B6').Style = 'Normal'
$ws1.Range('C6').Value = '南昌·CM04动漫游戏博览会'
$ws1.Range('D6').Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws1.Range('E6').Value = '2024.11.16 09:00-11.17 17:00'
$ws1.Range('F6').Value = 1486
$ws1.Range('G6').Value = 65
$ws1.Range('H6').Value = 'https://show.bilibili.com/platform/detail.html?id=92378'
$ws1.Range('I6').Value = '//i2.hdslb.com/bfs/openplatform/202409/N57Jfogr1725381095803.jpeg'
# row 7
$ws1.Range('B7').NumberFormat = '@'
$ws1.Range('B7').Value = '2024-11-30'
$ws1.Range('B7').Style = 'Normal'
$ws1.Range('C7').Value = '南昌·岁酉山河·炎国明日方舟同人ONLY'
$ws1.Range('D7').Value = '民德路411号 东方豪景花园酒店(民德路店)'
$ws1.Range('E7').Value = '2024.11.30 09:00-11.30 17:00'
$ws1.Range('F7').Value = 35
$ws1.Range('G7').Value = 68
$ws1.Range('H7').Value = 'https://show.bilibili.com/platform/detail.html?id=93050'
$ws1.Range('I7').Value = '//i2.hdslb.com/bfs/openplatform/202409/IBvdHJ1G1726720682507.png'
# row 8
$ws1.Range('B8').NumberFormat = '@'
$ws1.Range('B8').Value = '2024-12-07'
$ws1.Range('B8').Style = 'Normal'
$ws1.Range('C8').Value = '南昌·云芽动漫音乐嘉年华'
$ws1.Range('D8').Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws1.Range('E8').Value = '2024.12.07 09:00-12.08 18:00'
$ws1.Range('F8').Value = 3079
$ws1.Range('G8').Value = 69
$ws1.Range('H8').Value = 'https://show.bilibili.com/platform/detail.html?id=92144'
$ws1.Range('I8').Value = '//i0.hdslb.com/bfs/openplatform/202409/2DwZA4qv1725706772865.png'
# row 9
$ws1.Range('B9').NumberFormat = '@'
$ws1.Range('B9').Value = '2024-12-08'
$ws1.Range('B9').Style = 'Normal'
$ws1.Range('C9').Value = '南昌·云芽动漫音乐嘉年华·封茗囧菌内场票'
$ws1.Range('D9').Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws1.Range('E9').Value = '2024.12.08 09:30-12.08 17:30'
$ws1.Range('F9').Value = 451
$ws1.Range('G9').Value = '已售罄'
$ws1.Range('H9').Value = 'https://show.bilibili.com/platform/detail.html?id=92134'
$ws1.Range('I9').Value = '//i0.hdslb.com/bfs/openplatform/202409/eeFHJb3W1725328994111.jpeg'
# row 10
$ws1.Range('B10').NumberFormat = '@'
$ws1.Range('B10').Value = '2025-01-01'
$ws1.Range('B10').Style = 'Normal'
$ws1.Range('C10').Value = '南昌·萌卡动漫展'
$ws1.Range('D10').Value = '八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆'
$ws1.Range('E10').Value = '2025.01.01 09:00-01.03 17:00'
$ws1.Range('F10').Value = 711
$ws1.Range('G10').Value = 70
$ws1.Range('H10').Value = 'https://show.bilibili.com/platform/detail.html?id=93031'
$ws1.Range('I10').Value = '//i2.hdslb.com/bfs/openplatform/202409/HTlK8fN21727112669248.jpeg'

# remove the now-superseded trailing row (old row 11)
$ws1.Range('A11').EntireRow.Delete()

# --- Sheet '全部类型' (All types) ---
$ws4 = $wb.Worksheets.Item('全部类型')
# row 2
$ws4.Range('B2').NumberFormat = '@'
$ws4.Range('B2').Value = '2024-10-06'
$ws4.Range('B2').Style = 'Normal'
$ws4.Range('C2').Value = '南昌·星辰动漫游戏展嘉年华'
$ws4.Range('D2').Value = '长寿路 锐成体育2＋1篮球公园'
$ws4.Range('E2').Value = '2024.10.06 10:00-10.06 17:00'
$ws4.Range('F2').Value = 165
$ws4.Range('G2').Value = 45
$ws4.Range('H2').Value = 'https://show.bilibili.com/platform/detail.html?id=91961'
$ws4.Range('I2').Value = '//i2.hdslb.com/bfs/openplatform/202409/A6PiZvto1725362505262.jpeg'
# row 3
$ws4.Range('B3').NumberFormat = '@'
$ws4.Range('B3').Value = '2024-10-19'
$ws4.Range('B3').Style = 'Normal'
$ws4.Range('C3').Value = '南昌·第一届哥布林动漫游戏展——开学季&贺中秋'
$ws4.Range('D3').Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws4.Range('E3').Value = '2024.10.19 10:00-10.20 18:00'
$ws4.Range('F3').Value = 648
$ws4.Range('G3').Value = 55
$ws4.Range('H3').Value = 'https://show.bilibili.com/platform/detail.html?id=89240'
$ws4.Range('I3').Value = '//i2.hdslb.com/bfs/openplatform/202409/3N3iIqVu1725270119618.jpeg'
# row 4
$ws4.Range('B4').NumberFormat = '@'
$ws4.Range('B4').Value = '2024-11-02'
$ws4.Range('B4').Style = 'Normal'
$ws4.Range('C4').Value = '南昌·花绒万兽秋镜派对'
$ws4.Range('D4').Value = '双港西大街899号 旭辉Cmall(南昌店)'
$ws4.Range('E4').Value = '2024.11.02 10:00-11.03 21:30'
$ws4.Range('F4').Value = 24
$ws4.Range('G4').Value = 168
$ws4.Range('H4').Value = 'https://show.bilibili.com/platform/detail.html?id=92859'
$ws4.Range('I4').Value = '//i2.hdslb.com/bfs/openplatform/202409/7hJL2m3F1727175584690.jpeg'
# row 5
$ws4.Range('B5').NumberFormat = '@'
$ws4.Range('B5').Value = '2024-11-06'
$ws4.Range('B5').Style = 'Normal'
$ws4.Range('C5').Value = '南昌·松井祐贵 2024《阳光之旅》指弹吉他音乐会'
$ws4.Range('D5').Value = '上海路543号520Park文创公园21号01区域 瓦肆VAS NANCHANG'
$ws4.Range('E5').Value = '2024.11.06 20:00-11.06 21:30'
$ws4.Range('F5').Value = 1
$ws4.Range('G5').Value = 220
$ws4.Range('H5').Value = 'https://show.bilibili.com/platform/detail.html?id=92765'
$ws4.Range('I5').Value = '//i1.hdslb.com/bfs/openplatform/202409/iUNLvHVz1727082732931.jpeg'
# row 6
$ws4.Range('B6').NumberFormat = '@'
$ws4.Range('B6').Value = '2024-11-16'
$ws4.Range('B6').Style = 'Normal'
$ws4.Range('C6').Value = '上饶·星河城市动漫文化节'
$ws4.Range('D6').Value = '春江北大道时光PARK内 博悦宴会艺术中心'
$ws4.Range('E6').Value = '2024.11.16 10:00-11.16 17:00'
$ws4.Range('F6').Value = 214
$ws4.Range('G6').Value = 55
$ws4.Range('H6').Value = 'https://show.bilibili.com/platform/detail.html?id=92572'
$ws4.Range('I6').Value = '//i2.hdslb.com/bfs/openplatform/202409/xp4jNVRG1727165677359.jpeg'
# row 7
$ws4.Range('B7').NumberFormat = '@'
$ws4.Range('B7').Value = '2024-11-16'
$ws4.Range('B7').Style = 'Normal'
$ws4.Range('C7').Value = '南昌·CM04动漫游戏博览会'
$ws4.Range('D7').Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws4.Range('E7').Value = '2024.11.16 09:00-11.17 17:00'
$ws4.Range('F7').Value = 1486
$ws4.Range('G7').Value = 65
$ws4.Range('H7').Value = 'https://show.bilibili.com/platform/detail.html?id=92378'
$ws4.Range('I7').Value = '//i2.hdslb.com/bfs/openplatform/202409/N57Jfogr1725381095803.jpeg'
# row 8
$ws4.Range('B8').NumberFormat = '@'
$ws4.Range('B8').Value = '2024-11-30'
$ws4.Range('B8').Style = 'Normal'
$ws4.Range('C8').Value = '南昌·岁酉山河·炎国明日方舟同人ONLY'
$ws4.Range('D8').Value = '民德路411号 东方豪景花园酒店(民德路店)'
$ws4.Range('E8').Value = '2024.11.30 09:00-11.30 17:00'
$ws4.Range('F8').Value = 35
$ws4.Range('G8').Value = 68
$ws4.Range('H8').Value = 'https://show.bilibili.com/platform/detail.html?id=93050'
$ws4.Range('I8').Value = '//i2.hdslb.com/bfs/openplatform/202409/IBvdHJ1G1726720682507.png'
# row 9
$ws4.Range('B9').NumberFormat = '@'
$ws4.Range('B9').Value = '2024-12-07'
$ws4.Range('B9').Style = 'Normal'
$ws4.Range('C9').Value = '南昌·云芽动漫音乐嘉年华'
$ws4.Range('D9').Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws4.Range('E9').Value = '2024.12.07 09:00-12.08 18:00'
$ws4.Range('F9').Value = 3079
$ws4.Range('G9').Value = 69
$ws4.Range('H9').Value = 'https://show.bilibili.com/platform/detail.html?id=92144'
$ws4.Range('I9').Value = '//i0.hdslb.com/bfs/openplatform/202409/2DwZA4qv1725706772865.png'
# row 10
$ws4.Range('B10').NumberFormat = '@'
$ws4.Range('B10').Value = '2024-12-08'
$ws4.Range('B10').Style = 'Normal'
$ws4.Range('C10').Value = '南昌·云芽动漫音乐嘉年华·封茗囧菌内场票'
$ws4.Range('D10').Value = '怀玉山大道1315号 南昌绿地国际博览中心'
$ws4.Range('E10').Value = '2024.12.08 09:30-12.08 17:30'
$ws4.Range('F10').Value = 451
$ws4.Range('G10').Value = '已售罄'
$ws4.Range('H10').Value = 'https://show.bilibili.com/platform/detail.html?id=92134'
$ws4.Range('I10').Value = '//i0.hdslb.com/bfs/openplatform/202409/eeFHJb3W1725328994111.jpeg'
# row 11
$ws4.Range('B11').NumberFormat = '@'
$ws4.Range('B11').Value = '2025-01-01'
$ws4.Range('B11').Style = 'Normal'
$ws4.Range('C11').Value = '南昌·萌卡动漫展'
$ws4.Range('D11').Value = '八一桥街道青山南路118号蓝海购物广场F1 蓝海展览馆'
$ws4.Range('E11').Value = '2025.01.01 09:00-01.03 17:00'
$ws4.Range('F11').Value = 711
$ws4.Range('G11').Value = 70
$ws4.Range('H11').Value = 'https://show.bilibili.com/platform/detail.html?id=93031'
$ws4.Range('I11').Value = '//i2.hdslb.com/bfs/openplatform/202409/HTlK8fN21727112669248.jpeg'

# remove the now-superseded trailing row (old row 12)
$ws4.Range('A12').EntireRow.Delete()
